# =====================================================================
# Edit: rename sheet2, insert new "Database" info rows + chart updates,
# add new sheet "2023-05 Compr. Table Linux" with updated benchmark data.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Rename existing "Compression Table 2023-01" sheet
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "2023-01 Compr. Table Exadata"

# ---------------------------------------------------------------
# 2) Insert rows to make room for the "Database" info block
#    - old row1 (title) shifts to row2 (insert blank row at 1)
#    - old row3.. (Tested access onward) shifts by +3 total
#      (2 more blank rows inserted right before it)
# ---------------------------------------------------------------
$ws2.Rows.Item(1).Insert()
$ws2.Range("A4:A5").EntireRow.Insert()

# ---------------------------------------------------------------
# 3) Populate the new "Database" row (row 3)
# ---------------------------------------------------------------
$ws2.Range("A3").Value2 = "Database"
$ws2.Range("B3").Value2 = "19.18 auf Exadata X6-2L High Capacity "

# ---------------------------------------------------------------
# 4) Update the chart's series formulas (sheet rename + row shift)
#    and reposition/resize the chart so it still spans the same
#    rows relative to the table (rows 23-48 -> 26-51, 0-indexed).
# ---------------------------------------------------------------
$co = $ws2.ChartObjects().Item(1)
$chart = $co.Chart
$newSheetName = "2023-01 Compr. Table Exadata"

for ($i = 1; $i -le 7; $i++) {
    $ser = $chart.SeriesCollection().Item($i)
    $row = 15 + $i
    $f = "=SERIES('" + $newSheetName + "'!`$A`$" + $row + ",'" + $newSheetName + "'!`$B`$14:`$Q`$15,'" + $newSheetName + "'!`$B`$" + $row + ":`$Q`$" + $row + "," + $i + ")"
    $ser.Formula = $f
}

$newTop = $ws2.Rows.Item(27).Top + 12.5
$newBottom = $ws2.Rows.Item(52).Top + 7.0
$co.Top = $newTop
$co.Height = $newBottom - $newTop

# ---------------------------------------------------------------
# 5) Select B4 on sheet2 (matches final workbook selection state)
# ---------------------------------------------------------------
$ws2.Range("B4").Select()

# ---------------------------------------------------------------
# 6) Add the new "2023-05 Compr. Table Linux" sheet at the end,
#    by copying the fully-updated "Exadata" sheet's data+formats
#    (but not its chart) and then editing the values.
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "2023-05 Compr. Table Linux"

# Copy column widths from sheet2
for ($c = 1; $c -le 17; $c++) {
    $ws3.Columns.Item($c).ColumnWidth = $ws2.Columns.Item($c).ColumnWidth
}
$ws3.Columns.Item(16).ColumnWidth = 11.33203125

# Copy values+styles+merges for the shared header block (rows 2-15)
$ws2.Range("A2:Q15").Copy($ws3.Range("A2"))

# Copy row heights for header block
for ($r = 2; $r -le 15; $r++) {
    $ws3.Rows.Item($r).RowHeight = $ws2.Rows.Item($r).RowHeight
}

# Update the version text for the Linux database row
$ws3.Range("B3").Value2 = "19.3 auf Linux 6 CPU DB in Docker"

Write-Output "done with structural edits"
